$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Multi-Utilities(18)'
$ws.Range("B2").Value = 0.6806635493997961
$ws.Range("A3").Value = 'Household Products(10)'
$ws.Range("B3").Value = 0.5551239585304982
$ws.Range("A4").Value = 'Electric Utilities(28)'
$ws.Range("B4").Value = 0.5028578579683908
$ws.Range("A5").Value = 'Gas Utilities(12)'
$ws.Range("B5").Value = 0.4996720732227125
$ws.Range("A6").Value = 'Containers & Packaging(12)'
$ws.Range("B6").Value = 0.4554211694031084
$ws.Range("A7").Value = 'Insurance(75)'
$ws.Range("B7").Value = 0.4434552474521457
$ws.Range("A8").Value = 'Road & Rail(22)'
$ws.Range("B8").Value = 0.4197651903675268
$ws.Range("A9").Value = 'Marine(15)'
$ws.Range("B9").Value = 0.3924305106678252
$ws.Range("A10").Value = 'Professional Services(35)'
$ws.Range("B10").Value = 0.3515205790689357
$ws.Range("A11").Value = 'Construction & Engineering(20)'
$ws.Range("B11").Value = 0.3499373854978178
$ws.Range("A12").Value = 'Machinery(85)'
$ws.Range("B12").Value = 0.3261202385101064
$ws.Range("A13").Value = 'Building Products(23)'
$ws.Range("B13").Value = 0.3222612338074511
$ws.Range("A14").Value = 'Specialty Retail(58)'
$ws.Range("B14").Value = 0.3199989557889833
$ws.Range("A15").Value = 'Auto Components(21)'
$ws.Range("B15").Value = 0.3195755395178446
$ws.Range("A16").Value = 'Airlines(14)'
$ws.Range("B16").Value = 0.3116040117463391
$ws.Range("A17").Value = 'Food & Staples Retailing(15)'
$ws.Range("B17").Value = 0.3095211343106541
$ws.Range("A18").Value = 'Chemicals(51)'
$ws.Range("B18").Value = 0.2764694756478722
$ws.Range("A19").Value = 'Consumer Finance(15)'
$ws.Range("B19").Value = 0.275584626156352
$ws.Range("A20").Value = 'Diversified Telecommunication Services(20)'
$ws.Range("B20").Value = 0.2704492994662344
$ws.Range("A21").Value = 'Diversified Consumer Services(17)'
$ws.Range("B21").Value = 0.2576480116098388
$ws.Range("A22").Value = 'Banks(246)'
$ws.Range("B22").Value = 0.25398121351853
$ws.Range("A23").Value = 'Textiles, Apparel & Luxury Goods(29)'
$ws.Range("B23").Value = 0.2460280766957026
$ws.Range("A24").Value = 'Aerospace & Defense(37)'
$ws.Range("B24").Value = 0.2409153622487695
$ws.Range("A25").Value = 'Metals & Mining(89)'
$ws.Range("B25").Value = 0.231530059686375
$ws.Range("A26").Value = 'Food Products(44)'
$ws.Range("B26").Value = 0.2264936990572523
$ws.Range("A27").Value = 'Media(42)'
$ws.Range("B27").Value = 0.223127485609525
$ws.Range("A28").Value = 'Commercial Services & Supplies(52)'
$ws.Range("B28").Value = 0.2230322642084072
$ws.Range("A29").Value = 'Hotels, Restaurants & Leisure(50)'
$ws.Range("B29").Value = 0.2212959646850693
$ws.Range("A30").Value = 'Electrical Equipment(28)'
$ws.Range("B30").Value = 0.2202056081988574
$ws.Range("A31").Value = 'Energy Equipment & Services(32)'
$ws.Range("B31").Value = 0.2151226296637223
$ws.Range("A32").Value = 'IT Services(52)'
$ws.Range("B32").Value = 0.2133866368881704
$ws.Range("A33").Value = 'Trading Companies & Distributors(25)'
$ws.Range("B33").Value = 0.2072231145600254
$ws.Range("A34").Value = 'Capital Markets(75)'
$ws.Range("B34").Value = 0.2058945325758229
$ws.Range("A35").Value = 'Household Durables(39)'
$ws.Range("B35").Value = 0.1881941681680455
$ws.Range("A36").Value = 'Semiconductors & Semiconductor Equipment(68)'
$ws.Range("B36").Value = 0.1830530634682372
$ws.Range("A37").Value = 'Health Care Providers & Services(46)'
$ws.Range("B37").Value = 0.1693941077824654
$ws.Range("A38").Value = 'Health Care Equipment & Supplies(83)'
$ws.Range("B38").Value = 0.153903120944962
$ws.Range("A39").Value = 'Communications Equipment(45)'
$ws.Range("B39").Value = 0.1408272132494436
$ws.Range("A40").Value = 'Software(66)'
$ws.Range("B40").Value = 0.1382521755303223
$ws.Range("A41").Value = 'Pharmaceuticals(48)'
$ws.Range("B41").Value = 0.1327046731336035
$ws.Range("A42").Value = 'Thrifts & Mortgage Finance(47)'
$ws.Range("B42").Value = 0.1299470991891688
$ws.Range("A43").Value = 'Biotechnology(126)'
$ws.Range("B43").Value = 0.1235769747809436
$ws.Range("A44").Value = 'Oil, Gas & Consumable Fuels(122)'
$ws.Range("B44").Value = 0.09839953425492326

# Remove now-unused trailing rows 45 and 46
$ws.Range("A45:B46").Clear()
